# Reorders the comma-separated Net Name / Net Label lists in the "Net Name"
# (col AB) and "Net Label" (col AC) columns of the BoM and DNF sheets.
# These are plain text cell values (not formulas); the underlying XLSX uses
# a shared-string table, but each cell must be (re)written individually
# because writing one cell does not retroactively repoint other cells that
# happened to reference the same shared-string entry.

$wb = $excel.ActiveWorkbook

$bom = $wb.Worksheets.Item("BoM")

$bom.Range("AB10").Value = "Net-(U1-UCAP),Earth"
$bom.Range("AC10").Value = "Net-(U1-UCAP),Earth"

$bom.Range("AB12").Value = "+5V,/RESET2"
$bom.Range("AC12").Value = "RESET2"

$bom.Range("AB13").Value = "Net-(D2-A),/RXLED"
$bom.Range("AC13").Value = "RXLED"

$bom.Range("AB15").Value = "+5V,/SCK2,/MISO2,/RESET2,/MOSI2,GND"
$bom.Range("AC15").Value = "MOSI2,GND"

$bom.Range("AB16").Value = "Net-(J3-Pin_5),Net-(J3-Pin_1),Net-(J3-Pin_4),Net-(J3-Pin_3),Net-(J3-Pin_2)"
$bom.Range("AC16").Value = "Net-(J3-Pin_5),Net-(J3-Pin_1),Net-(J3-Pin_4),Net-(J3-Pin_3),Net-(J3-Pin_2)"

$bom.Range("AB17").Value = "Net-(J6-Pin_3),Net-(J6-Pin_2),Net-(J6-Pin_6),Net-(J6-Pin_1),Net-(J6-Pin_4),Net-(J6-Pin_5)"
$bom.Range("AC17").Value = "Net-(J6-Pin_3),Net-(J6-Pin_2),Net-(J6-Pin_6),Net-(J6-Pin_1),Net-(J6-Pin_4),Net-(J6-Pin_5)"

$bom.Range("AB20").Value = "+5V,/RESET2"
$bom.Range("AC20").Value = "RESET2"

$bom.Range("AB21").Value = "Net-(J3-Pin_5),VBUS,Net-(J3-Pin_4),unconnected-(U1-PB0-Pad14),/SCK2,Net-(U1-UCAP),Net-(J6-Pin_4),/TXLED,/DTR,Net-(U1-D-),/MISO2,Earth,Net-(J3-Pin_1),Net-(J6-Pin_3),Net-(J6-Pin_2),Net-(J3-Pin_3),Net-(U1-D+),Net-(J6-Pin_5),Net-(J3-Pin_2),Net-(J4-Pin_2),Net-(U1-XTAL1),Net-(J4-Pin_3),Net-(J6-Pin_6),Net-(J4-Pin_1),Net-(U1-PC0{slash}XTAL2),+5V,Net-(J4-Pin_4),/RESET2,/MOSI2,/RXLED,GND"
$bom.Range("AC21").Value = "RXLED,GND"

$dnf = $wb.Worksheets.Item("DNF")

$dnf.Range("AB9").Value = "GND,Net-(U1-XTAL1)"
$dnf.Range("AC9").Value = "GND,Net-(U1-XTAL1)"

$dnf.Range("AB10").Value = "Net-(J2-VBUS),VBUS"
$dnf.Range("AC10").Value = "Net-(J2-VBUS),VBUS"

$dnf.Range("AB12").Value = "Net-(J4-Pin_1),Net-(J4-Pin_3),Net-(J4-Pin_2),Net-(J4-Pin_4)"
$dnf.Range("AC12").Value = "Net-(J4-Pin_1),Net-(J4-Pin_3),Net-(J4-Pin_2),Net-(J4-Pin_4)"

$dnf.Range("AB13").Value = "Net-(J2-VBUS),Net-(J2-Shield),Net-(J2-D+),Earth,Net-(J2-D-)"
$dnf.Range("AC13").Value = "Net-(J2-VBUS),Net-(J2-Shield),Net-(J2-D+),Earth,Net-(J2-D-)"

$dnf.Range("AB14").Value = "Net-(U1-D+),Net-(J2-D+)"
$dnf.Range("AC14").Value = "Net-(U1-D+),Net-(J2-D+)"

$dnf.Range("AB15").Value = "Net-(U1-PC0{slash}XTAL2),Net-(U1-XTAL1)"
$dnf.Range("AC15").Value = "Net-(U1-PC0{slash}XTAL2),Net-(U1-XTAL1)"

$dnf.Range("AB16").Value = "Net-(J2-D+),Net-(J2-Shield)"
$dnf.Range("AC16").Value = "Net-(J2-D+),Net-(J2-Shield)"

$dnf.Range("AB17").Value = "Net-(U1-PC0{slash}XTAL2),Net-(U1-XTAL1)"
$dnf.Range("AC17").Value = "Net-(U1-PC0{slash}XTAL2),Net-(U1-XTAL1)"
